$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = '60.047.93'
$ws.Cells.Item(2, 5).Value = '  +2.43%  '

# Row 3
$ws.Cells.Item(3, 4).Value = '3.193.68'
$ws.Cells.Item(3, 5).Value = '  +1.31%  '

# Row 4
$ws.Cells.Item(4, 5).Value = '  +0.00%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '536.47'
$ws.Cells.Item(5, 5).Value = '  +1.00%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '145.17'
$ws.Cells.Item(6, 5).Value = '  +3.91%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = '0.531'
$ws.Cells.Item(8, 5).Value = '  -1.09%  '

# Row 9
$ws.Cells.Item(9, 5).Value = '  +0.15%  '

# Row 10
$ws.Cells.Item(10, 5).Value = '  +2.26%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '0.432'
$ws.Cells.Item(11, 5).Value = '  -0.65%  '

# Row 12
$ws.Cells.Item(12, 4).Value = '3.742.83'
$ws.Cells.Item(12, 5).Value = '  +1.30%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = '0.138'
$ws.Cells.Item(13, 5).Value = '  -2.33%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '25.82'
$ws.Cells.Item(14, 5).Value = '  -0.16%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '0.0000172'
$ws.Cells.Item(15, 5).Value = '  +0.97%  '

# Row 16
$ws.Cells.Item(16, 4).Value = '60.035.78'
$ws.Cells.Item(16, 5).Value = '  +2.33%  '

# Row 17
$ws.Cells.Item(17, 4).Value = '3.194.29'
$ws.Cells.Item(17, 5).Value = '  +0.98%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '6.26'
$ws.Cells.Item(18, 5).Value = '  +0.09%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '13.28'
$ws.Cells.Item(19, 5).Value = '  +2.16%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '8.21'
$ws.Cells.Item(20, 5).Value = '  +0.72%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '370.01'
$ws.Cells.Item(21, 5).Value = '  -0.35%  '

# Row 22
$ws.Cells.Item(22, 5).Value = '  +0.01%  '

# Row 23
$ws.Cells.Item(23, 5).Value = '  -0.48%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '69.41'
$ws.Cells.Item(24, 5).Value = '  -0.30%  '

# Row 25
$ws.Cells.Item(25, 5).Value = '  +1.28%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '8.57'
$ws.Cells.Item(26, 5).Value = '  +4.07%  '

# Row 27
$ws.Cells.Item(27, 5).Value = '  -0.99%  '

# Row 28
$ws.Cells.Item(28, 4).Value = '0.0₃0874'
$ws.Cells.Item(28, 5).Value = '  +1.31%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '22.49'
$ws.Cells.Item(29, 5).Value = '  +1.96%  '

# Row 30
$ws.Cells.Item(30, 5).Value = '  +0.65%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '6.11'
$ws.Cells.Item(31, 5).Value = '  +0.08%  '

# Row 32
$ws.Cells.Item(32, 5).Value = '  +2.57%  '

# Row 33
$ws.Cells.Item(33, 5).Value = '  +2.75%  '

# Row 34
$ws.Cells.Item(34, 5).Value = '  +4.05%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '156.60'
$ws.Cells.Item(35, 5).Value = '  -1.25%  '

# Row 36
$ws.Cells.Item(36, 5).Value = '  +1.89%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '26.57'
$ws.Cells.Item(37, 5).Value = '  +5.87%  '

# Row 38
$ws.Cells.Item(38, 4).Value = '2.819.56'
$ws.Cells.Item(38, 5).Value = '  +7.04%  '

# Row 39
$ws.Cells.Item(39, 5).Value = '  +2.96%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '0.0311'
$ws.Cells.Item(40, 5).Value = '  +8.68%  '

# Row 41
$ws.Cells.Item(41, 5).Value = '  -0.41%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '4.24'
$ws.Cells.Item(42, 5).Value = '  -0.22%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '39.59'
$ws.Cells.Item(43, 5).Value = '  +1.58%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '0.719'
$ws.Cells.Item(44, 5).Value = '  +1.69%  '

# Row 45
$ws.Cells.Item(45, 2).Value = 'RenzoRestakedETH'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Cells.Item(45, 4).Value = '3.235.57'
$ws.Cells.Item(45, 5).Value = '  +1.27%  '

# Row 46
$ws.Cells.Item(46, 2).Value = 'Stellar'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D46").NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '0.104'
$ws.Cells.Item(46, 5).Value = '  +0.60%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '0.984'
$ws.Cells.Item(47, 5).Value = '  +0.41%  '

# Row 48
$ws.Cells.Item(48, 5).Value = '  -0.80%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '20.65'
$ws.Cells.Item(49, 5).Value = '  +2.06%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '0.793'
$ws.Cells.Item(50, 5).Value = '  +4.99%  '

# Row 51
$ws.Cells.Item(51, 5).Value = '  +0.00%  '
